# Autogenerated on Thu Mar 26 2015 18:06:15 GMT+0000 (Coordinated Universal Time)
#
# Reworks the "Source" citation block at the bottom of the Canada Summary
# sheet (rows 52-59 -> 52-60):
#   - A53's text ("Industry Canada - Key Small Business Statistics...")
#     moves down to A54, and the hyperlink that used to live on A54 is
#     dropped (the URL becomes a plain text line further down instead).
#   - A53 becomes blank.
#   - A new row is inserted after A55 to hold the URL as plain text (A56).
#   - The "INDCA" abbreviation (old A58) shifts down to A59.
#   - The long citation paragraph (old A59) is replaced by a second,
#     duplicate "INDCA" line at A60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row right after the current A55 ("" / source style).
# Excel pushes the old rows 56-59 (blank gap, INDCA, long citation) down to
# 57-60, which lines up A58->A59 and A59->A60 exactly as the diff wants.
$ws.Range("A56").EntireRow.Insert()

# The old A54 hyperlink (pointing at the Industry Canada URL) no longer
# applies once the cell's text changes, so drop it first.
$ws.Range("A54").Hyperlinks.Delete()

# A53 loses its text -- it becomes the blank "source"-styled cell.
$ws.Range("A53").Value2 = ""
$ws.Range("A53").Font.Italic = $true
$ws.Range("A53").Font.Underline = $false

# A54 now holds the description text that used to be on A53, styled like
# the other plain "source" citation lines (italic, no hyperlink look).
$ws.Range("A54").Value2 = "Industry Canada - Key Small Business Statistics - August 2013"
$ws.Range("A54").Font.Italic = $true
$ws.Range("A54").Font.Underline = $false

# The newly inserted row holds the URL as plain (non-hyperlinked) text.
$ws.Range("A56").Value2 = "http://www.ic.gc.ca/eic/site/061.nsf/eng/h_02800.html"
$ws.Range("A56").Font.Italic = $true
$ws.Range("A56").Font.Underline = $false

# A59 ("INDCA") keeps the bold "title" look of the old A58 it replaced --
# row-insert already carried the value/style down, nothing else to do.
$ws.Range("A59").Font.Bold = $true

# A60 replaces the long citation paragraph with a second "INDCA" line,
# keeping the italic "source" look of the old A59 it replaced.
$ws.Range("A60").Value2 = "INDCA"
$ws.Range("A60").Font.Italic = $true
$ws.Range("A60").Font.Underline = $false
